# Actualización automática 2025-09-26 08:30:09
# Update "CUMPLIMIENTO MENSUAL" sheet with new VENTA (sales) figures for
# rows 3 (240X80 PORCELANATO) and 12 (PORCELANATO), recomputing the
# dependent "POR CUMPLIR" and "CUMPLIMIENTO" columns, as well as the
# TOTAL row (15).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- Row 3: 240X80 PORCELANATO ---
$ws.Range("D3").Value = 3967.49
$ws.Range("E3").Value = -1385.07619324963
$ws.Range("F3").Value = 1.536349437734988

# --- Row 12: PORCELANATO ---
$ws.Range("D12").Value = 5019.5
$ws.Range("E12").Value = 12655.8486842162
$ws.Range("F12").Value = 0.2839830822959851

# --- Row 15: TOTAL ---
$ws.Range("D15").Value = 15298.37
$ws.Range("E15").Value = 16409.38990313501
$ws.Range("F15").Value = 0.4824803154412501
